# Config.xlsx update: "Ajustado arquivo de configuração"
# Rebrands the REFramework config from "ProcessABCQueue"/"Framework" to a
# "Cadastro" project, adds a Shared queue folder, and a new email/Gmail
# credential block to the Settings sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

$xlCenter = -4108
$xlLeft   = -4131

# --- Row 1 header: "Value" column header becomes bold + left-aligned ---
$ws.Range("B1").HorizontalAlignment = $xlLeft

# --- Row 2: OrchestratorQueueName value -> Cadastro ---
$ws.Range("B2").Value = "Cadastro"
$ws.Range("B2").HorizontalAlignment = $xlCenter

# --- Row 3: OrchestratorQueueFolder value -> Shared ---
$ws.Range("B3").Value = "Shared"
$ws.Range("B3").HorizontalAlignment = $xlCenter
$ws.Rows.Item(3).RowHeight = 45

# --- Row 5: logF_BusinessProcessName value -> Cadastro ---
$ws.Range("B5").Value = "Cadastro"
$ws.Range("B5").HorizontalAlignment = $xlCenter
$ws.Rows.Item(5).RowHeight = 30

# --- Row 8: new CadastroUrl setting ---
$ws.Range("A8").Value = "CadastroUrl"
$ws.Range("B8").Value = "https://forms.office.com/r/WmF5T2PXBj"
$ws.Range("B8").HorizontalAlignment = $xlCenter

# --- Row 10: new EmailCredential setting ---
$ws.Range("A10").Value = "EmailCredential"
$ws.Range("B10").Value = "EmailCredential"
$ws.Range("B10").HorizontalAlignment = $xlCenter

# --- Row 12: new GmailPort setting ---
$ws.Range("A12").Value = "GmailPort"
$ws.Range("B12").Value = 465
$ws.Range("B12").HorizontalAlignment = $xlCenter

# --- Row 14: new GmailServer setting ---
$ws.Range("A14").Value = "GmailServer"
$ws.Range("B14").Value = "imap.gmail.com"
$ws.Range("B14").HorizontalAlignment = $xlCenter

# --- Move the saved selection cursor the way the author left it ---
$ws.Range("B15").Select()

$wb.Save()
